$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$zwsp = [char]0x200B

# "conservative" table - O row (row 2): update e_geom, e_norm, Intensity
$ws.Range("C2").Value = "7.89e-10$zwsp"
$ws.Range("D2").Value = "2.88E-6$zwsp"
$ws.Range("E2").Value = "1.81E9$zwsp"

# "25 ns conservative" table - O row (row 22): update e_geom, e_norm, Intensity
$ws.Range("C22").Value = "7.91e-10$zwsp"
$ws.Range("D22").Value = "2.89E-6$zwsp"
$ws.Range("E22").Value = "1.8E9$zwsp"

# Update the selection/view state
$ws.Range("C23").Select()
